$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = "r567"
$ws.Range("B6").Value = "fred"
$ws.Range("C6").Value = "very scary"
$ws.Range("D6").Value = "2025-09-30 20:29:00"
